# Applies the "gh-pages output generated at 456a3b4" update:
#   - refresh "想去人数" (F column) counters on sheet 展览 (1) and 全部类型 (4)
#   - refresh F12/F13 counters on sheet 演出 (2)
#   - insert a new event row (2024-08-11, Marcin Patrzalek concert) as the
#     new row 16 on sheet 演出 (2), pushing the existing row 16 down to row 17

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibition) - sheet index 1
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F2").Value  = 7608
$wsExpo.Range("F3").Value  = 7608
$wsExpo.Range("F5").Value  = 7796
$wsExpo.Range("F9").Value  = 6508
$wsExpo.Range("F10").Value = 3337
$wsExpo.Range("F12").Value = 3691
$wsExpo.Range("F15").Value = 37
$wsExpo.Range("F17").Value = 44
$wsExpo.Range("F22").Value = 3779
$wsExpo.Range("F24").Value = 365
$wsExpo.Range("F26").Value = 279
$wsExpo.Range("F27").Value = 1436
$wsExpo.Range("F30").Value = 2721
$wsExpo.Range("F31").Value = 1742
$wsExpo.Range("F35").Value = 3571
$wsExpo.Range("F36").Value = 276
$wsExpo.Range("F40").Value = 523
$wsExpo.Range("F41").Value = 1381
$wsExpo.Range("F42").Value = 240
$wsExpo.Range("F44").Value = 628

# ---------------------------------------------------------------------
# Sheet "演出" (Performance) - sheet index 2
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F12").Value = 57
$wsShow.Range("F13").Value = 85

# Insert a brand-new row 16 (shifts the old row 16 -> row 17, and bumps
# its serial number in column A from 15 to 16).
$wsShow.Rows.Item(16).Insert()

$wsShow.Range("A16").Value = 15
# Force text so the ISO-looking date string is stored verbatim instead of
# being auto-converted to a date serial number (matches the rest of column B,
# which is plain text in the source data).
$wsShow.Range("B16").NumberFormat = "@"
$wsShow.Range("B16").Value = "2024-08-11"
$wsShow.Range("C16").Value = "北京·Marcin Patrzalek 2024 《原声之龙》指弹吉他音乐会"
$wsShow.Range("D16").Value = "西坝河南里2号香河园地区文化中心 多维剧场"
$wsShow.Range("E16").Value = "2024.08.11 20:00-08.11 21:30"
$wsShow.Range("F16").Value = 0
$wsShow.Range("G16").Value = 380
$wsShow.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=86309"
$wsShow.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202405/MEqm9GHU1716777275477.jpeg"

# Match the style Excel uses for the bold/bordered serial-number column
# (same formatting as every other cell in column A).
$wsShow.Range("A16").Font.Bold = $true
$wsShow.Range("A16").HorizontalAlignment = -4108
$wsShow.Range("A16").VerticalAlignment = -4160
$wsShow.Range("A16").Borders.LineStyle = 1

# Row 17 (the event that used to be row 16) keeps its own data untouched
# except its running counter in column A, which advances from 15 to 16.
$wsShow.Range("A17").Value = 16

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) - sheet index 4
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F5").Value  = 7608
$wsAll.Range("F6").Value  = 7608
$wsAll.Range("F8").Value  = 7796
$wsAll.Range("F11").Value = 6508
$wsAll.Range("F12").Value = 3337
$wsAll.Range("F14").Value = 3691
$wsAll.Range("F17").Value = 37
$wsAll.Range("F19").Value = 44
$wsAll.Range("F24").Value = 3779
$wsAll.Range("F28").Value = 57
$wsAll.Range("F29").Value = 365
$wsAll.Range("F31").Value = 279
$wsAll.Range("F32").Value = 1436
$wsAll.Range("F35").Value = 2721
$wsAll.Range("F36").Value = 1742
$wsAll.Range("F39").Value = 85
$wsAll.Range("F40").Value = 3571
$wsAll.Range("F41").Value = 276
$wsAll.Range("F42").Value = 272
$wsAll.Range("F45").Value = 523
$wsAll.Range("F46").Value = 1381
$wsAll.Range("F47").Value = 240
$wsAll.Range("F50").Value = 628
